# #195 adding changes for production readiness
#
# - "Sheet1" (India venues) becomes the first tab (was second) and stays
#   the active/selected tab.
# - "Sheet2" (UAE venues) becomes the second tab (was first) and is no
#   longer the active tab.
# - "Feroz Shah Kotla" (DELHI) is replaced by "Arun Jaitley Stadium".
# - A new row for "Narendra Modi Stadium" / "Ahmedabad" / "INDIA" is added.

$wb = $excel.ActiveWorkbook

$indiaSheet = $wb.Worksheets.Item("Sheet1")
$uaeSheet   = $wb.Worksheets.Item("Sheet2")

# --- Data edits (done while sheet references are still fresh) ---

# New row 9: Narendra Modi Stadium / Ahmedabad / INDIA
$indiaSheet.Cells.Item(9, 1).Value = "Narendra Modi Stadium"
$indiaSheet.Cells.Item(9, 2).Value = "Ahmedabad"
$indiaSheet.Cells.Item(9, 3).Value = "INDIA"

# Match the look of the other "away" venue rows (Eden Gardens, Sawai
# Mansingh, Chinnaswamy, Rajiv Gandhi stadiums) for the new row.
$indiaSheet.Range("A8").Copy()
$indiaSheet.Range("A9").PasteSpecial(-4122)
$indiaSheet.Rows.Item(9).RowHeight = 20

# Feroz Shah Kotla -> Arun Jaitley Stadium (row 2, DELHI)
$indiaSheet.Cells.Item(2, 1).Value = "Arun Jaitley Stadium"

# --- Selections on each sheet before reordering ---
[void]$uaeSheet.Range("B30").Select()
[void]$indiaSheet.Range("C7").Select()

# --- Reorder tabs: India sheet ("Sheet1") becomes the first tab ---
$indiaSheet.Move($uaeSheet)

# Re-resolve by name (sheet refs track tab position in this host) and make
# sure the India sheet is the active tab with its intended selection.
$indiaSheet = $wb.Worksheets.Item("Sheet1")
[void]$indiaSheet.Range("C7").Select()

Write-Output "Reordered sheets and updated IPL venue data."
